# Adding ship class and ship class output
# - wrap-text style applied to the header row (A1:I1)
# - "Did Sensed State Match Truth State?" (column I) formula simplified
#   from IF(AND(Hn=1, An=Fn), 1, 0) to IF(Bn=Gn, 1, 0) for every data row
# - several Sensor column values (C, E, F) updated for rows describing
#   the new ship-class scenarios

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Apply wrap-text formatting to the header row, A1:I1
$ws.Range("A1:I1").WrapText = $true

# 2. Update the "Did Sensed State Match Truth State?" formula (column I)
#    for every data row (rows 2 through 122)
for ($r = 2; $r -le 122; $r++) {
    $ws.Range("I$r").Formula = "=IF(B$r=G$r, 1, 0)"
}

# 3. Update the Sensor values that changed with the new ship-class data
$ws.Range("E5:E62").Value = 0
$ws.Range("E65:E122").Value = 0

$ws.Range("C11:C15").Value = 0
$ws.Range("C71:C75").Value = 0

$ws.Range("F12:F62").Value = 1
$ws.Range("F72:F122").Value = 1
